$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ItemProto")

# Write the name strings first (column B) so they land at shared-string
# indices 68/69, then the icon path strings (column C) at 70/71 - this
# matches the order new <si> entries were appended to sharedStrings.xml.
$ws.Cells.Item(18, 2).Value = "精炼符"
$ws.Cells.Item(19, 2).Value = "超级精炼符"
$ws.Cells.Item(18, 3).Value = "icon/item/item_4111.png"
$ws.Cells.Item(19, 3).Value = "icon/item/item_4112.png"

# New row 18: item 4111 - 精炼符
$ws.Cells.Item(18, 1).Value = 4111
$ws.Cells.Item(18, 5).Value = 5
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 12).Value = 0

# New row 19: item 4112 - 超级精炼符
$ws.Cells.Item(19, 1).Value = 4112
$ws.Cells.Item(19, 5).Value = 5
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 12).Value = 0

# D column stays blank but keeps the same "s=2" formatting as the rest of
# the table - copy formats only (values untouched) from row 17.
$ws.Cells.Item(17, 2).Copy()
$ws.Cells.Item(18, 2).PasteSpecial(-4122)
$ws.Cells.Item(17, 3).Copy()
$ws.Cells.Item(18, 3).PasteSpecial(-4122)
$ws.Cells.Item(17, 4).Copy()
$ws.Cells.Item(18, 4).PasteSpecial(-4122)

$ws.Cells.Item(17, 2).Copy()
$ws.Cells.Item(19, 2).PasteSpecial(-4122)
$ws.Cells.Item(17, 3).Copy()
$ws.Cells.Item(19, 3).PasteSpecial(-4122)
$ws.Cells.Item(17, 4).Copy()
$ws.Cells.Item(19, 4).PasteSpecial(-4122)

$excel.CutCopyMode = $false

$ws.Range("C23").Select()
